$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The status column (F) values for these non-conformities were still "En proceso";
# update them to "Cerrada" to reflect that they have been closed.
$ws.Range("F8").Value = "Cerrada"
$ws.Range("F9").Value = "Cerrada"
$ws.Range("F11").Value = "Cerrada"

# Restore the view/selection state left after making the edits.
$ws.Range("F12").Select()
